# download models & update results
# Apply the experiment-results update to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 12 (J12/K12 "经常折返" / "不是稳赢，经常折返") is removed outright -
# that content is being relocated down into row 14 (J14/K14 below).
# ---------------------------------------------------------------------
$ws.Range("J12:K12").Clear()

# ---------------------------------------------------------------------
# New annotation / result cells, written in the same order the author
# typed them so newly-introduced shared strings line up.
# ---------------------------------------------------------------------

# Run-dir note for the SD row.
$ws.Range("O13").Value = "map*_use_step_dist"

# Run-dir note for the SD-3F row (was mis-named before).
$ws.Range("O14").Value = "map*_use_step_dist_[frames]3"

# Row 14 (SD-3F) per-map results.
$ws.Range("H3").Copy()
$ws.Range("M14").PasteSpecial(-4122)
$ws.Range("M14").Value = "只能过个弯"

$ws.Range("G3").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("D14").Value = "可以，不流畅"

$ws.Range("G3").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$ws.Range("E14").Value = "可以，不流畅"

$ws.Range("L6").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("F14").Value = "经常走回头路，经常不过"

$ws.Range("L6").Copy()
$ws.Range("G14").PasteSpecial(-4122)
$ws.Range("G14").Value = "经常走回头路，偶尔能过"

$ws.Range("H3").Copy()
$ws.Range("H14").PasteSpecial(-4122)
$ws.Range("H14").Value = "走的远一点，但是会卡角落"

$ws.Range("H3").Copy()
$ws.Range("I14").PasteSpecial(-4122)
$ws.Range("I14").Value = "不行"

$ws.Range("H3").Copy()
$ws.Range("L14").PasteSpecial(-4122)
$ws.Range("L14").Value = "只能过两三个弯"

# Row 13 (SD) gains new annotation cells.
$ws.Range("G3").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("C13").Value = "很标准沿中间走"

$ws.Range("G3").Copy()
$ws.Range("J13").PasteSpecial(-4122)
$ws.Range("J13").Value = "偶尔折返"

$ws.Range("L6").Copy()
$ws.Range("K13").PasteSpecial(-4122)
$ws.Range("K13").Value = "不是稳赢，经常折返+卡墙"

# Row 14 cells reusing text that used to live in row 12.
$ws.Range("G3").Copy()
$ws.Range("J14").PasteSpecial(-4122)
$ws.Range("J14").Value = "经常折返"

$ws.Range("L6").Copy()
$ws.Range("K14").PasteSpecial(-4122)
$ws.Range("K14").Value = "不是稳赢，经常折返"

# ---------------------------------------------------------------------
# View-state touch-ups (zoom + active selection).
# ---------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 125
$ws.Range("I12").Select() | Out-Null
